$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh crypto price (D) and 1h volume change (E) columns with latest scraped data.
# Some "Price" values are plain numbers; force them to be stored as text (matching
# the source feed formatting, e.g. multi-dot big numbers like "26.174.55") so Excel
# does not silently reinterpret them as numeric values.

$ws.Range("D2").Value = '26.174.55'
$ws.Range("E2").Value = '  -0.07%  '
$ws.Range("D3").Value = '1.670.72'
$ws.Range("E3").Value = '  -0.66%  '
$ws.Range("E4").Value = '  -0.34%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.46'
$ws.Range("E5").Value = '  -2.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5221'
$ws.Range("E6").Value = '  -0.45%  '
$ws.Range("E7").Value = '  -0.36%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2627'
$ws.Range("E8").Value = '  -2.39%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06331'
$ws.Range("E9").Value = '  -0.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.27'
$ws.Range("E10").Value = '  -1.00%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07554'
$ws.Range("E11").Value = '  -1.11%  '
$ws.Range("D12").Value = '1.679.47'
$ws.Range("E12").Value = '  -0.36%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.444'
$ws.Range("E13").Value = '  -1.62%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5507'
$ws.Range("E14").Value = '  -4.20%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '66.61'
$ws.Range("E15").Value = '  +0.95%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008001'
$ws.Range("E16").Value = '  -3.62%  '
$ws.Range("D17").Value = '26.189.33'
$ws.Range("E17").Value = '  -0.27%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.003'
$ws.Range("E18").Value = '  -0.36%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.756'
$ws.Range("E19").Value = '  -2.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '187.32'
$ws.Range("E20").Value = '  -0.84%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.224'
$ws.Range("E22").Value = '  -0.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.004'
$ws.Range("E23").Value = '  -0.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '149.82'
$ws.Range("E24").Value = '  +0.94%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1246'
$ws.Range("E25").Value = '  -1.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.513'
$ws.Range("E26").Value = '  -3.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.88'
$ws.Range("E27").Value = '  +0.87%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06327'
$ws.Range("E28").Value = '  +0.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.354'
$ws.Range("E29").Value = '  -1.69%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.279'
$ws.Range("E30").Value = '  -2.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.521'
$ws.Range("E31").Value = '  -1.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.418'
$ws.Range("E32").Value = '  -4.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.646'
$ws.Range("E33").Value = '  -1.94%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.003'
$ws.Range("E34").Value = '  -1.97%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6033'
$ws.Range("E35").Value = '  -1.29%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.407'
$ws.Range("E36").Value = '  -0.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.757'
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.146'
$ws.Range("E38").Value = '  -0.77%  '
$ws.Range("D39").Value = '1.110.52'
$ws.Range("E39").Value = '  +1.42%  '
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8653'
$ws.Range("E41").Value = '  -2.91%  '
$ws.Range("E42").Value = '  -0.57%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.35'
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("D44").Value = '1.825.38'
$ws.Range("E44").Value = '  -0.39%  '
$ws.Range("E45").Value = '  -1.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '55.57'
$ws.Range("E46").Value = '  -2.95%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.003'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.056'
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05235'
$ws.Range("E49").Value = '  -0.84%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4242'
$ws.Range("E50").Value = '  -0.92%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.931'
$ws.Range("E51").Value = '  -1.41%  '
